# Insert a new data row at row 338 (shifting existing rows 338:402 down to 339:403)
# and populate it with a new price observation, matching the rest of the
# "Femacal de La Calera - Arándano (blue)" dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("338:338").Insert()

$ws.Range("A338").Value = 3
$ws.Range("B338").Value = "Femacal de La Calera"
$ws.Range("C338").Value = "Coquimbo"
$ws.Range("D338").Value = 45244
$ws.Range("E338").Value = 5
$ws.Range("F338").Value = "Fruta"
$ws.Range("G338").Value = 100101
$ws.Range("H338").Value = "Berries"
$ws.Range("I338").Value = 100101001
$ws.Range("J338").Value = "Arándano (blue)"
$ws.Range("K338").Value = "Sin especificar"
$ws.Range("L338").Value = "Primera"
$ws.Range("M338").Value = 54
$ws.Range("N338").Value = 6000
$ws.Range("O338").Value = 6000
$ws.Range("P338").Value = 6000
$ws.Range("Q338").Value = "`$/bandeja 2 kilos"
$ws.Range("R338").Value = "Provincia de Quillota"
$ws.Range("S338").Value = 3000
$ws.Range("T338").Value = 2
